# gh-pages data refresh: updated "想去人数" (interested-count) figures.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Sheets.Item("展览")
$wsExpo.Range("F3").Value = 48
$wsExpo.Range("F4").Value = 2280
$wsExpo.Range("F6").Value = 503

# Sheet "全部类型" (All types) mirrors the same rows
$wsAll = $wb.Sheets.Item("全部类型")
$wsAll.Range("F5").Value = 48
$wsAll.Range("F6").Value = 2280
$wsAll.Range("F8").Value = 503
